$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "Logs" sheet: append the new row (row 10) with the fresh ticket.
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Cells.Item(10, 1).Value = "Demo inplannen"
$logs.Cells.Item(10, 2).Value = "klantenservice@testbedrijf123.nl"
$logs.Cells.Item(10, 3).Value = "Kun je vrijdag om 11:00 een demo inplannen bij Van Dijk?"
$logs.Cells.Item(10, 4).Value = "Planning / Afspraak"
$logs.Cells.Item(10, 5).Value = "Bedankt, we hebben dit doorgestuurd naar planning@testbedrijf123.nl."
$logs.Cells.Item(10, 6).Value = "2025-08-14 20:37:26"
$logs.Cells.Item(10, 7).Value = "Nee"
$logs.Cells.Item(10, 8).Value = "Ja"
$logs.Cells.Item(10, 9).Value = "Nee"
$logs.Cells.Item(10, 10).Value = "Nee"

# Extend the conditional-formatting blocks (D/G/H/I/J) so the new row is
# covered too, mirroring how Excel grows these ranges when rows are added.
$logs.Range("D2:D9").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D10"))
$logs.Range("G2:G9").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G10"))
$logs.Range("H2:H9").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H10"))
$logs.Range("I2:I9").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I10"))
$logs.Range("J2:J9").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J10"))

# ---------------------------------------------------------------------
# 2) "Dashboard" sheet: add the new category count row (row 3).
# ---------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Cells.Item(3, 1).Value = "Planning / Afspraak"
$dash.Cells.Item(3, 2).Value = 1

# ---------------------------------------------------------------------
# 3) Chart on the Dashboard sheet: extend the category/value series refs
#    so the new Dashboard row is plotted too.
# ---------------------------------------------------------------------
$chart = $dash.ChartObjects(1).Chart
$ser = $chart.SeriesCollection(1)
$ser.XValues = '=''Dashboard''!$A$2:$A$3'
$ser.Values = '=''Dashboard''!$B$2:$B$3'
